$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'290.65"
$ws.Range("E2").Value = "'-5.95%"

# Row 3
$ws.Range("D3").Value = "'39.78"
$ws.Range("E3").Value = "'-3.17%"

# Row 4
$ws.Range("D4").Value = "'5.016"
$ws.Range("E4").Value = "'-3.20%"

# Row 5
$ws.Range("D5").Value = "'0.07348"
$ws.Range("E5").Value = "'-4.04%"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.553"
$ws.Range("E6").Value = "'-11.19%"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9186"
$ws.Range("E7").Value = "'0.33%"

# Row 8
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1190"
$ws.Range("E8").Value = "'-4.13%"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1726"
$ws.Range("E9").Value = "'-4.79%"

# Row 10
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.08696"
$ws.Range("E10").Value = "'-4.58%"

# Row 11
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04172"
$ws.Range("E11").Value = "'0.34%"

# Row 12
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.1054"
$ws.Range("E12").Value = "'0.43%"

# Row 13
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001276"
$ws.Range("E13").Value = "'1.17%"

# Row 14
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005785"
$ws.Range("E14").Value = "'0.11%"

# Row 15
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.406"
$ws.Range("E15").Value = "'1.83%"

# Row 16
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'4.287"
$ws.Range("E16").Value = "'-0.29%"

# Row 17
$ws.Range("D17").Value = "'2.397"
$ws.Range("E17").Value = "'-1.16%"

# Row 18
$ws.Range("E18").Value = "'-0.69%"

# Row 19
$ws.Range("D19").Value = "'7.589"
$ws.Range("E19").Value = "'2.87%"

# Row 20
$ws.Range("E20").Value = "'-0.16%"

# Row 22
$ws.Range("D22").Value = "'0.03849"
$ws.Range("E22").Value = "'-4.16%"

# Row 23
$ws.Range("E23").Value = "'0.63%"

# Row 24
$ws.Range("D24").Value = "'0.003886"
$ws.Range("E24").Value = "'-5.13%"

# Row 25
$ws.Range("D25").Value = "'0.0001283"
$ws.Range("E25").Value = "'-1.67%"

# Row 26
$ws.Range("D26").Value = "'0.0003731"

# Row 38
$ws.Range("D38").Value = "'0.02338"
$ws.Range("E38").Value = "'-7.30%"

# Row 39
$ws.Range("D39").Value = "'0.05020"
$ws.Range("E39").Value = "'-5.12%"

# Row 40
$ws.Range("D40").Value = "'0.007684"
$ws.Range("E40").Value = "'-2.21%"

# Row 41
$ws.Range("E41").Value = "'172.28%"

# Row 42
$ws.Range("E42").Value = "'-2.93%"

# Row 43
$ws.Range("D43").Value = "'0.007381"
$ws.Range("E43").Value = "'11.08%"

# Row 44
$ws.Range("D44").Value = "'0.007718"
$ws.Range("E44").Value = "'-4.82%"

# Row 45
$ws.Range("D45").Value = "'0.3173"
$ws.Range("E45").Value = "'3.47%"

# Row 46
$ws.Range("D46").Value = "'0.00006518"
$ws.Range("E46").Value = "'-4.17%"

# Row 47
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'-0.03%"

# Row 48
$ws.Range("E48").Value = "'8.93%"

# Row 49
$ws.Range("D49").Value = "'0.004209"
$ws.Range("E49").Value = "'35.63%"

# Row 50
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'-0.03%"

# Row 51
$ws.Range("E51").Value = "'-0.03%"
